$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its text formatting so numeric-looking
# strings like "26.553.42" are not coerced into Excel numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.553.42"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.812.41"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D6").Value = "305.95"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").Value = "0.4550"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "0.3596"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").Value = "46.38"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("D10").Value = "0.07110"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "0.8933"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("D12").Value = "0.07719"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "19.33"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "1.800.16"
$ws.Range("D15").Value = "5.253"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "6.291"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "86.65"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "0.000008548"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "26.583.49"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "14.15"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "4.960"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("D26").Value = "151.75"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "17.80"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "2.018"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").Value = "112.43"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "4.831"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "0.08723"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "3.124"
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").Value = "0.7391"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "2.730"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.435"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "1.111"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "0.01935"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("D39").Value = "2.915"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "0.05072"
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").Value = "0.5079"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").Value = "6.783"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").Value = "0.1506"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").Value = "8.016"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "0.4687"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "9.990"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "99.14"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").Value = "1.566"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").Value = "0.06002"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "63.61"
$ws.Range("E51").Value = "  -1.44%  "
